# Update the metadata for the "municipio-nombre" column (column J) so that
# it is now described using the newly curated dimension info instead of
# being treated as a measure.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "sdmx-dimension:refArea"
$ws.Range("J3").Value = "dim"
$ws.Range("J4").Value = "URI-Municipio"
